$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New header for column F
$ws.Range("F1").Value = "Image File"

# Fill F2:F25 with the "Q##.PNG" formula and G2:G25 with "ok"
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("F$r").Formula = '="Q"&IF(LEN(A' + $r + ')=1,"0","")&A' + $r + '&".PNG"'
    $ws.Range("G$r").Value = "ok"
}

# Match the style of column F to the rest of the table (centered alignment)
$ws.Range("F1:F25").HorizontalAlignment = -4108

# Update the selected cell to reflect the new active cell
$ws.Range("G26").Select()
